# Update column G ("K") values on the active sheet to reflect the
# regenerated save_data (switch from Strike# to K, recalculated values).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = 0
    3 = 0
    4 = 0
    5 = 1
    6 = 2
    7 = 4
    8 = 1
    9 = 3
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
